# pMeHg_OutflowR.xlsx edit script
# Implements:
#  - rename worksheet "4_pMeHg" -> "4_pMeHg_cens"
#  - change column D header "pMeHg" -> "PpMeHg" and add a new column E
#    (header "PpMeHg" / unit row "2s") to flag censored ("<") values
#  - update the reporting/detection limit for rows 4 & 5 (0.05 -> 0.1) and
#    mark those rows as censored ("<") in the new column E
#  - change the highlighted ("non-detect") fill color used for the
#    censored-value cells from a plain gold to the Gold/Accent4 theme color

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the worksheet ---------------------------------------------
$ws.Name = "4_pMeHg_cens"

# --- header rows --------------------------------------------------------
# Row 1: column headers. D stays "PpMeHg" (renamed from "pMeHg") and the
# new column E repeats the same header.
$ws.Range("D1").Value = "PpMeHg"
$ws.Range("E1").Value = "PpMeHg"

# Row 2: parameter/code row. Existing A2:D2 content is unchanged logically
# (10d / 8s / 6n / 6n); only the new E2 cell is added.
$ws.Range("E2").Value = "2s"

# --- updated censored values (rows 4 & 5) --------------------------------
# Reporting limit changed from 0.05 to 0.1, flagged as "<" (less-than /
# non-detect) in the new column E. The cells keep their existing
# highlighted ("non-detect") style automatically.
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = "<"

$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = "<"

# --- recolor the "non-detect" highlight fill -----------------------------
# Previously a flat gold (RGB FFC000); now uses the theme's Accent4 (gold)
# color family instead of a hard-coded RGB value.
$ws.Range("D4").Interior.ThemeColor = 8
$ws.Range("D5").Interior.ThemeColor = 8
